$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controller")
$ws.Activate()

# Add the three new rows of test data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "CHROME"
$ws.Range("C4").Value = "Test_03"
$ws.Range("D4").Value = "Yes"
$ws.Range("E4").Value = "PASS"
$ws.Range("F4").Value = "DEV1"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "CHROME"
$ws.Range("C5").Value = "Test_04"
$ws.Range("D5").Value = "Yes"
$ws.Range("E5").Value = "PASS"
$ws.Range("F5").Value = "DEV1"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "CHROME"
$ws.Range("C6").Value = "Test_05"
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "FAIL"
$ws.Range("F6").Value = "DEV1"

# Copy style formatting from row 3 to the new rows (A:G) to match existing format
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false | Out-Null

# Update selection to match target state
$ws.Range("A3:A6").Select() | Out-Null
